$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 122 and 123: swap the match data (columns F:V) between the two
#    rows. Columns A:E (index/pais/torneio/temporada/data_partida) stay as
#    they are.
# ---------------------------------------------------------------------------

# NOTE: this runtime's `.Value` *getter* does not surface the underlying
# scalar (it echoes member metadata instead); `.Value2` reads correctly
# (and, unlike `.Text`, preserves numeric typing instead of formatted
# text), so use it for all reads below. `.Value` remains fine as a setter.

# Capture row 122's current F:V values before overwriting anything.
$f122 = $ws.Cells.Item(122, 6).Value2
$g122 = $ws.Cells.Item(122, 7).Value2
$h122 = $ws.Cells.Item(122, 8).Value2
$i122 = $ws.Cells.Item(122, 9).Value2
$j122 = $ws.Cells.Item(122, 10).Value2
$k122 = $ws.Cells.Item(122, 11).Value2
$l122 = $ws.Cells.Item(122, 12).Value2
$m122 = $ws.Cells.Item(122, 13).Value2
$n122 = $ws.Cells.Item(122, 14).Value2
$o122 = $ws.Cells.Item(122, 15).Value2
$p122 = $ws.Cells.Item(122, 16).Value2
$q122 = $ws.Cells.Item(122, 17).Value2
$r122 = $ws.Cells.Item(122, 18).Value2
$s122 = $ws.Cells.Item(122, 19).Value2
$t122 = $ws.Cells.Item(122, 20).Value2
$u122 = $ws.Cells.Item(122, 21).Value2
$v122 = $ws.Cells.Item(122, 22).Value2

# Capture row 123's current F:V values.
$f123 = $ws.Cells.Item(123, 6).Value2
$g123 = $ws.Cells.Item(123, 7).Value2
$h123 = $ws.Cells.Item(123, 8).Value2
$i123 = $ws.Cells.Item(123, 9).Value2
$j123 = $ws.Cells.Item(123, 10).Value2
$k123 = $ws.Cells.Item(123, 11).Value2
$l123 = $ws.Cells.Item(123, 12).Value2
$m123 = $ws.Cells.Item(123, 13).Value2
$n123 = $ws.Cells.Item(123, 14).Value2
$o123 = $ws.Cells.Item(123, 15).Value2
$p123 = $ws.Cells.Item(123, 16).Value2
$q123 = $ws.Cells.Item(123, 17).Value2
$r123 = $ws.Cells.Item(123, 18).Value2
$s123 = $ws.Cells.Item(123, 19).Value2
$t123 = $ws.Cells.Item(123, 20).Value2
$u123 = $ws.Cells.Item(123, 21).Value2
$v123 = $ws.Cells.Item(123, 22).Value2

# Write the old row-123 data into row 122.
$ws.Cells.Item(122, 6).Value = $f123
$ws.Cells.Item(122, 7).Value = $g123
$ws.Cells.Item(122, 8).Value = $h123
$ws.Cells.Item(122, 9).Value = $i123
$ws.Cells.Item(122, 10).Value = $j123
$ws.Cells.Item(122, 11).Value = $k123
$ws.Cells.Item(122, 12).Value = $l123
$ws.Cells.Item(122, 13).Value = $m123
$ws.Cells.Item(122, 14).Value = $n123
$ws.Cells.Item(122, 15).Value = $o123
$ws.Cells.Item(122, 16).Value = $p123
$ws.Cells.Item(122, 17).Value = $q123
$ws.Cells.Item(122, 18).Value = $r123
$ws.Cells.Item(122, 19).Value = $s123
$ws.Cells.Item(122, 20).Value = $t123
$ws.Cells.Item(122, 21).Value = $u123
$ws.Cells.Item(122, 22).Value = $v123

# Write the old row-122 data into row 123.
$ws.Cells.Item(123, 6).Value = $f122
$ws.Cells.Item(123, 7).Value = $g122
$ws.Cells.Item(123, 8).Value = $h122
$ws.Cells.Item(123, 9).Value = $i122
$ws.Cells.Item(123, 10).Value = $j122
$ws.Cells.Item(123, 11).Value = $k122
$ws.Cells.Item(123, 12).Value = $l122
$ws.Cells.Item(123, 13).Value = $m122
$ws.Cells.Item(123, 14).Value = $n122
$ws.Cells.Item(123, 15).Value = $o122
$ws.Cells.Item(123, 16).Value = $p122
$ws.Cells.Item(123, 17).Value = $q122
$ws.Cells.Item(123, 18).Value = $r122
$ws.Cells.Item(123, 19).Value = $s122
$ws.Cells.Item(123, 20).Value = $t122
$ws.Cells.Item(123, 21).Value = $u122
$ws.Cells.Item(123, 22).Value = $v122

# ---------------------------------------------------------------------------
# 2) Append two new match rows (150 and 151) at the bottom of the sheet.
#    Column A (Indice) and column E (data_partida) carry the same formatting
#    as the rest of the table (bold/bordered index, date-formatted match
#    date), so clone that formatting from an existing row before writing
#    the values.
# ---------------------------------------------------------------------------

$ws.Cells.Item(149, 1).Copy() | Out-Null
$ws.Cells.Item(150, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(151, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(149, 5).Copy() | Out-Null
$ws.Cells.Item(150, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(151, 5).PasteSpecial(-4122) | Out-Null

# Row 150
$ws.Cells.Item(150, 1).Value = 149
$ws.Cells.Item(150, 2).Value = "serbia"
$ws.Cells.Item(150, 3).Value = "super-liga"
$ws.Cells.Item(150, 4).Value = "2023-2024"
$ws.Cells.Item(150, 5).Value = 45282.66666666666
$ws.Cells.Item(150, 6).Value = "Javor"
$ws.Cells.Item(150, 7).Value = 1
$ws.Cells.Item(150, 8).Value = "IMT Novi Beograd"
$ws.Cells.Item(150, 9).Value = 1
$ws.Cells.Item(150, 10).Value = 2.25
$ws.Cells.Item(150, 11).Value = "26/09/2023 16:13"
$ws.Cells.Item(150, 12).Value = 2.38
$ws.Cells.Item(150, 13).Value = "22/12/2023 15:59"
$ws.Cells.Item(150, 14).Value = 3.1
$ws.Cells.Item(150, 15).Value = "26/09/2023 16:13"
$ws.Cells.Item(150, 16).Value = 3.18
$ws.Cells.Item(150, 17).Value = "22/12/2023 15:58"
$ws.Cells.Item(150, 18).Value = 2.88
$ws.Cells.Item(150, 19).Value = "26/09/2023 16:13"
$ws.Cells.Item(150, 20).Value = 2.97
$ws.Cells.Item(150, 21).Value = "22/12/2023 15:59"
$ws.Cells.Item(150, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/javor-imt-novi-beograd/MmCHhtpU/"

# Row 151
$ws.Cells.Item(151, 1).Value = 150
$ws.Cells.Item(151, 2).Value = "serbia"
$ws.Cells.Item(151, 3).Value = "super-liga"
$ws.Cells.Item(151, 4).Value = "2023-2024"
$ws.Cells.Item(151, 5).Value = 45282.70833333334
$ws.Cells.Item(151, 6).Value = "TSC"
$ws.Cells.Item(151, 7).Value = 1
$ws.Cells.Item(151, 8).Value = "Napredak"
$ws.Cells.Item(151, 9).Value = 1
$ws.Cells.Item(151, 10).Value = 1.33
$ws.Cells.Item(151, 11).Value = "25/09/2023 13:12"
$ws.Cells.Item(151, 12).Value = 1.34
$ws.Cells.Item(151, 13).Value = "22/12/2023 16:44"
$ws.Cells.Item(151, 14).Value = 4.58
$ws.Cells.Item(151, 15).Value = "25/09/2023 13:12"
$ws.Cells.Item(151, 16).Value = 4.86
$ws.Cells.Item(151, 17).Value = "22/12/2023 16:45"
$ws.Cells.Item(151, 18).Value = 7.33
$ws.Cells.Item(151, 19).Value = "25/09/2023 13:12"
$ws.Cells.Item(151, 20).Value = 8.57
$ws.Cells.Item(151, 21).Value = "22/12/2023 16:45"
$ws.Cells.Item(151, 22).Value = "https://www.betexplorer.com/football/serbia/super-liga/tsc-backa-topola-napredak/IeF0dOp5/"

# ---------------------------------------------------------------------------
# 3) Keep the sheet's declared dimension in sync with the new extent.
# ---------------------------------------------------------------------------
$ws.Range("A1:V151").Select() | Out-Null
